$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) changes - force text format to preserve exact string representation
# (values like "0.05339" or "26.539.66" would otherwise be parsed as numbers by Excel)
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D35","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.539.66"
$ws.Range("D3").Value = "1.729.99"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D5").Value = "245.59"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D7").Value = "0.4810"
$ws.Range("D8").Value = "0.2672"
$ws.Range("D9").Value = "0.06221"
$ws.Range("D10").Value = "1.727.37"
$ws.Range("D11").Value = "0.07164"
$ws.Range("D12").Value = "15.71"
$ws.Range("D13").Value = "0.6176"
$ws.Range("D14").Value = "4.537"
$ws.Range("D15").Value = "77.22"
$ws.Range("D17").Value = "26.543.31"
$ws.Range("D18").Value = "0.9995"
$ws.Range("D19").Value = "0.000006952"
$ws.Range("D20").Value = "11.67"
$ws.Range("D21").Value = "1.949.48"
$ws.Range("D22").Value = "4.534"
$ws.Range("D23").Value = "8.938"
$ws.Range("D24").Value = "5.291"
$ws.Range("D25").Value = "136.75"
$ws.Range("D26").Value = "15.36"
$ws.Range("D27").Value = "1.797"
$ws.Range("D28").Value = "1.404"
$ws.Range("D29").Value = "107.07"
$ws.Range("D30").Value = "3.991"
$ws.Range("D31").Value = "0.08040"
$ws.Range("D32").Value = "3.718"
$ws.Range("D33").Value = "0.04581"
$ws.Range("D35").Value = "2.616"
$ws.Range("D37").Value = "0.9962"
$ws.Range("D38").Value = "0.9245"
$ws.Range("D39").Value = "2.091"
$ws.Range("D40").Value = "2.408"
$ws.Range("D41").Value = "104.80"
$ws.Range("D42").Value = "1.007"
$ws.Range("D43").Value = "0.01505"
$ws.Range("D44").Value = "5.593"
$ws.Range("D45").Value = "0.3902"
$ws.Range("D46").Value = "6.967"
$ws.Range("D47").Value = "0.1184"
$ws.Range("D48").Value = "0.05339"
$ws.Range("D49").Value = "31.03"
$ws.Range("D50").Value = "7.826"
$ws.Range("D51").Value = "1.266"

# Column E (Volume/1h change %) changes - plain text assignment (never numeric-looking)
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +2.89%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +1.85%  "
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("E12").Value = "  +2.65%  "
$ws.Range("E13").Value = "  +4.60%  "
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("E23").Value = "  +2.17%  "
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("E27").Value = "  +2.50%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  +3.91%  "
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  +3.01%  "
$ws.Range("E37").Value = "  +1.77%  "
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("E39").Value = "  +9.78%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -7.78%  "
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("E43").Value = "  +1.96%  "
$ws.Range("E44").Value = "  +4.69%  "
$ws.Range("E45").Value = "  +2.49%  "
$ws.Range("E46").Value = "  +10.77%  "
$ws.Range("E47").Value = "  +1.53%  "
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("E50").Value = "  +2.19%  "
